# Add a new item row ("Machined Metal" / 가공된 금속) to Sheet1 of the Item Sheet workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 7

# Fill the new row's cells in column order so that newly introduced shared
# strings are appended to the shared-string table in the same left-to-right
# order used in the source file (Name, Name_eng, ItemType, Description, ...).
$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "가공된 금속"
$ws.Cells.Item($newRow, 3).Value = "MachinedMatal"
$ws.Cells.Item($newRow, 4).Value = "Key"
$ws.Cells.Item($newRow, 5).Value = "특이한 생김새로 가공된 금속이다. 어딘가에 쓰일지도?"
$ws.Cells.Item($newRow, 6).Value = 1
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = -1

# Move the active selection, matching the saved state in the file
$ws.Range("G8").Select()
